$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OR stunting by compfeeding")
$ws.Range("A1").Value = "test"
